$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look numeric need to be forced to Text
# so Excel doesn't silently coerce them (dropping trailing zeros, etc.).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '63.921.47'
$ws.Range('E2').Value = '  +1.42%  '
Set-TextValue 'D3' '3.303.88'
$ws.Range('E3').Value = '  +5.70%  '
Set-TextValue 'D5' '598.81'
$ws.Range('E5').Value = '  +0.60%  '
Set-TextValue 'D6' '142.99'
$ws.Range('E6').Value = '  +4.50%  '
$ws.Range('E7').Value = '  +0.00%  '
Set-TextValue 'D8' '3.298.53'
$ws.Range('E8').Value = '  +5.76%  '
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('E10').Value = '  +2.20%  '
Set-TextValue 'D11' '5.44'
$ws.Range('E11').Value = '  +3.25%  '
$ws.Range('E12').Value = '  +2.27%  '
$ws.Range('E13').Value = '  -0.44%  '
Set-TextValue 'D14' '34.81'
$ws.Range('E14').Value = '  +1.20%  '
Set-TextValue 'D15' '3.851.46'
$ws.Range('E15').Value = '  +5.83%  '
$ws.Range('E16').Value = '  +1.00%  '
Set-TextValue 'D17' '3.309.32'
$ws.Range('E17').Value = '  +5.84%  '
Set-TextValue 'D18' '63.998.85'
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('E19').Value = '  +2.03%  '
Set-TextValue 'D20' '481.21'
$ws.Range('E20').Value = '  +0.97%  '
Set-TextValue 'D21' '14.28'
$ws.Range('E21').Value = '  +0.46%  '
Set-TextValue 'D22' '0.743'
$ws.Range('E22').Value = '  +6.04%  '
$ws.Range('E23').Value = '  +4.24%  '
Set-TextValue 'D24' '13.48'
$ws.Range('E24').Value = '  +3.23%  '
$ws.Range('E25').Value = '  -3.51%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  +2.36%  '
$ws.Range('E28').Value = '  +1.69%  '
$ws.Range('E29').Value = '  -0.05%  '
Set-TextValue 'D30' '8.25'
$ws.Range('E30').Value = '  +3.37%  '
$ws.Range('E31').Value = '  +2.04%  '
Set-TextValue 'D32' '28.49'
$ws.Range('E32').Value = '  +4.86%  '
Set-TextValue 'D33' '0.107'
$ws.Range('E33').Value = '  -1.61%  '
$ws.Range('E34').Value = '  +0.66%  '
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('E36').Value = '  +2.57%  '
Set-TextValue 'D37' '53.33'
$ws.Range('E37').Value = '  +2.72%  '
$ws.Range('E38').Value = '  +3.71%  '
$ws.Range('E39').Value = '  +2.28%  '
Set-TextValue 'D40' '431.80'
$ws.Range('E40').Value = '  +1.86%  '
Set-TextValue 'D41' '3.015.53'
$ws.Range('E41').Value = '  +4.42%  '
Set-TextValue 'D42' '8.43'
$ws.Range('E42').Value = '  +1.74%  '
$ws.Range('E43').Value = '  +3.02%  '
$ws.Range('E44').Value = '  -6.36%  '
$ws.Range('E45').Value = '  +1.71%  '
Set-TextValue 'D46' '2.22'
$ws.Range('E46').Value = '  +4.36%  '
$ws.Range('E47').Value = '  +1.74%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D49' '0.115'
$ws.Range('E49').Value = '  +1.53%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D50' '2.33'
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('E51').Value = '  +13.90%  '
